# Automatische test-sync: 2025-06-20 11:00:50
# Append the new "Offerte / Prijsaanvraag" mail-log entry to the Logs sheet
# and roll the Dashboard summary + chart ranges forward to include it.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 8 ------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D8").Value = "Offerte / Prijsaanvraag"
$logs.Range("F8").Value = "2025-06-20 11:00:12"
$logs.Range("G8").Value = "Nee"

# Extend the conditional formatting ranges so they keep covering the full
# Categorie (D) and Beantwoord (G) columns now that row 8 exists.
$catFormats = $logs.Range("D2:D7").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D8"))
}

$answeredFormats = $logs.Range("G2:G7").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G8"))
}

# --- Dashboard sheet: append summary row 6 -----------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A6").Value = "Offerte / Prijsaanvraag"
$dashboard.Range("B6").Value = 1

# Roll the bar chart's category/value series ranges forward to include row 6.
$chartObj = $dashboard.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$6,Dashboard!`$B`$2:`$B`$6,1)"

Write-Output "Logs row 8 and Dashboard row 6 added; ranges extended."
